$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 252
$ws.Range("F5").Value = 2416
$ws.Range("F7").Value = 139
$ws.Range("F10").Value = 4829
$ws.Range("F12").Value = 861
$ws.Range("F13").Value = 79
$ws.Range("F14").Value = 1374
$ws.Range("F15").Value = 1327
$ws.Range("F16").Value = 525
$ws.Range("F17").Value = 6743
$ws.Range("F18").Value = 373
$ws.Range("F21").Value = 4513
$ws.Range("F22").Value = 355
$ws.Range("F25").Value = 2167
$ws.Range("F27").Value = 399
$ws.Range("F29").Value = 155
$ws.Range("F34").Value = 1238
$ws.Range("F35").Value = 1946
$ws.Range("F36").Value = 185
$ws.Range("F38").Value = 184
$ws.Range("F39").Value = 1294
$ws.Range("F40").Value = 574
$ws.Range("F42").Value = 51
$ws.Range("F43").Value = 1030
$ws.Range("F44").Value = 1317
$ws.Range("F46").Value = 85
$ws.Range("F48").Value = 51

# Sheet 2: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 420
$ws.Range("F5").Value = 442
$ws.Range("F15").Value = 145
$ws.Range("F19").Value = 122
$ws.Range("F35").Value = 20

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 1607
$ws.Range("G8").Value = "不可售"
$ws.Range("F9").Value = 1178
$ws.Range("F11").Value = 1664
$ws.Range("F12").Value = 1996
$ws.Range("F13").Value = 457
$ws.Range("F14").Value = 358

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1607
$ws.Range("F5").Value = 2416
$ws.Range("F6").Value = 1178
$ws.Range("F8").Value = 1664
$ws.Range("F10").Value = 1996
$ws.Range("F11").Value = 4829
$ws.Range("F12").Value = 442
$ws.Range("F14").Value = 861
$ws.Range("F15").Value = 79
$ws.Range("F17").Value = 1374
$ws.Range("F18").Value = 1327
$ws.Range("F19").Value = 525
$ws.Range("F20").Value = 6743
$ws.Range("F21").Value = 373
$ws.Range("F22").Value = 358
$ws.Range("F24").Value = 4513
$ws.Range("F25").Value = 355
$ws.Range("F27").Value = 399
$ws.Range("F29").Value = 155
$ws.Range("F35").Value = 1238
$ws.Range("F36").Value = 1946
$ws.Range("F37").Value = 185
$ws.Range("F39").Value = 122
$ws.Range("F40").Value = 184
$ws.Range("F41").Value = 1294
$ws.Range("F43").Value = 574
$ws.Range("F46").Value = 1317
$ws.Range("F47").Value = 85
$ws.Range("F49").Value = 51
